# Weekly GitHub Actions data refresh for the cryptos list (coinranking.com
# scrape). For every coin row this updates the "Price" (column D) and
# "Volume(1h)" (column E) figures to the latest snapshot. A few coins
# changed relative rank since the previous run, so for those rows the
# "Coin" (B) and "Link" (C) values are also rewritten so each row keeps
# describing one consistent coin.
#
# The Price column holds numeric-looking text (e.g. "1.012", "30.563.22")
# that must stay plain text (same as the original inline strings), but
# assigning such a string straight to .Value makes Excel auto-convert it
# to a real number (and, for tiny values, into scientific notation).
# Briefly forcing a text NumberFormat while assigning, then restoring the
# default "Normal" style, keeps the values as text without leaving any
# extra formatting on the cell.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Cells.Item(2, 4) "30.563.22"
$ws.Cells.Item(2, 5).Value = "  -1.47%  "
Set-TextValue $ws.Cells.Item(3, 4) "1.947.79"
$ws.Cells.Item(3, 5).Value = "  -1.34%  "
Set-TextValue $ws.Cells.Item(4, 4) "1.012"
$ws.Cells.Item(4, 5).Value = "  +1.56%  "
Set-TextValue $ws.Cells.Item(5, 4) "249.40"
$ws.Cells.Item(5, 5).Value = "  -1.42%  "
Set-TextValue $ws.Cells.Item(6, 4) "0.6900"
$ws.Cells.Item(6, 5).Value = "  -14.27%  "
Set-TextValue $ws.Cells.Item(7, 4) "1.012"
$ws.Cells.Item(7, 5).Value = "  +1.90%  "
Set-TextValue $ws.Cells.Item(8, 4) "0.3254"
$ws.Cells.Item(8, 5).Value = "  -4.91%  "
Set-TextValue $ws.Cells.Item(9, 4) "26.48"
$ws.Cells.Item(9, 5).Value = "  +3.77%  "
Set-TextValue $ws.Cells.Item(10, 4) "0.06805"
$ws.Cells.Item(10, 5).Value = "  -1.28%  "
Set-TextValue $ws.Cells.Item(11, 4) "0.7985"
$ws.Cells.Item(11, 5).Value = "  -6.31%  "
Set-TextValue $ws.Cells.Item(12, 4) "0.08014"
$ws.Cells.Item(12, 5).Value = "  -1.90%  "
Set-TextValue $ws.Cells.Item(13, 4) "1.961.26"
$ws.Cells.Item(13, 5).Value = "  -0.08%  "
Set-TextValue $ws.Cells.Item(14, 4) "5.424"
$ws.Cells.Item(14, 5).Value = "  -1.21%  "
Set-TextValue $ws.Cells.Item(15, 4) "94.54"
$ws.Cells.Item(15, 5).Value = "  -7.11%  "
Set-TextValue $ws.Cells.Item(16, 4) "14.56"
$ws.Cells.Item(16, 5).Value = "  +4.97%  "
Set-TextValue $ws.Cells.Item(17, 4) "262.51"
$ws.Cells.Item(17, 5).Value = "  -6.05%  "
Set-TextValue $ws.Cells.Item(18, 4) "30.616.19"
$ws.Cells.Item(18, 5).Value = "  -1.13%  "
Set-TextValue $ws.Cells.Item(19, 4) "5.911"
$ws.Cells.Item(19, 5).Value = "  +4.33%  "
Set-TextValue $ws.Cells.Item(20, 4) "0.000007862"
$ws.Cells.Item(20, 5).Value = "  +0.28%  "
Set-TextValue $ws.Cells.Item(21, 4) "2.223.18"
$ws.Cells.Item(21, 5).Value = "  +0.15%  "
Set-TextValue $ws.Cells.Item(22, 4) "1.010"
$ws.Cells.Item(22, 5).Value = "  +2.19%  "
Set-TextValue $ws.Cells.Item(23, 4) "1.013"
$ws.Cells.Item(23, 5).Value = "  +1.72%  "
Set-TextValue $ws.Cells.Item(24, 4) "6.878"
$ws.Cells.Item(24, 5).Value = "  +2.14%  "
Set-TextValue $ws.Cells.Item(25, 4) "9.710"
$ws.Cells.Item(25, 5).Value = "  +1.28%  "
Set-TextValue $ws.Cells.Item(26, 4) "158.90"
$ws.Cells.Item(26, 5).Value = "  -3.54%  "
Set-TextValue $ws.Cells.Item(27, 4) "18.90"
$ws.Cells.Item(27, 5).Value = "  -3.07%  "
Set-TextValue $ws.Cells.Item(28, 4) "2.274"
$ws.Cells.Item(28, 5).Value = "  +4.07%  "
Set-TextValue $ws.Cells.Item(29, 4) "0.1295"
$ws.Cells.Item(29, 5).Value = "  -22.04%  "
Set-TextValue $ws.Cells.Item(30, 4) "1.384"
$ws.Cells.Item(30, 5).Value = "  +2.34%  "
Set-TextValue $ws.Cells.Item(31, 4) "1.568"
$ws.Cells.Item(31, 5).Value = "  +0.54%  "
Set-TextValue $ws.Cells.Item(32, 4) "4.436"
$ws.Cells.Item(32, 5).Value = "  -1.95%  "
Set-TextValue $ws.Cells.Item(33, 4) "4.243"
$ws.Cells.Item(33, 5).Value = "  -1.78%  "
Set-TextValue $ws.Cells.Item(34, 4) "0.05118"
$ws.Cells.Item(34, 5).Value = "  +0.15%  "
Set-TextValue $ws.Cells.Item(35, 4) "1.205"
$ws.Cells.Item(35, 5).Value = "  -1.01%  "
Set-TextValue $ws.Cells.Item(36, 4) "0.7516"
$ws.Cells.Item(36, 5).Value = "  +1.92%  "
Set-TextValue $ws.Cells.Item(37, 4) "2.745"
$ws.Cells.Item(37, 5).Value = "  -0.42%  "
Set-TextValue $ws.Cells.Item(38, 4) "0.01950"
$ws.Cells.Item(38, 5).Value = "  -1.53%  "
Set-TextValue $ws.Cells.Item(39, 4) "2.837"
$ws.Cells.Item(39, 5).Value = "  -1.84%  "
Set-TextValue $ws.Cells.Item(40, 4) "80.34"
$ws.Cells.Item(40, 5).Value = "  +2.83%  "
Set-TextValue $ws.Cells.Item(41, 4) "6.601"
$ws.Cells.Item(41, 5).Value = "  +0.50%  "
$ws.Cells.Item(42, 2).Value = "RenderToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Cells.Item(42, 4) "2.052"
$ws.Cells.Item(42, 5).Value = "  -0.68%  "
$ws.Cells.Item(43, 2).Value = "TheSandbox"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Cells.Item(43, 4) "0.4446"
$ws.Cells.Item(43, 5).Value = "  -4.35%  "
$ws.Cells.Item(44, 2).Value = "PaxDollar"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Cells.Item(44, 4) "1.012"
$ws.Cells.Item(44, 5).Value = "  +1.91%  "
$ws.Cells.Item(45, 2).Value = "TrustWalletToken"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Cells.Item(45, 4) "0.8440"
$ws.Cells.Item(45, 5).Value = "  -0.21%  "
Set-TextValue $ws.Cells.Item(46, 4) "101.76"
$ws.Cells.Item(46, 5).Value = "  -2.46%  "
Set-TextValue $ws.Cells.Item(47, 4) "9.736"
$ws.Cells.Item(47, 5).Value = "  -2.21%  "
Set-TextValue $ws.Cells.Item(48, 4) "7.327"
$ws.Cells.Item(48, 5).Value = "  -1.63%  "
Set-TextValue $ws.Cells.Item(49, 4) "36.20"
$ws.Cells.Item(49, 5).Value = "  +0.06%  "
$ws.Cells.Item(50, 2).Value = "NEARProtocol"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Cells.Item(50, 4) "1.501"
$ws.Cells.Item(50, 5).Value = "  +5.02%  "
$ws.Cells.Item(51, 2).Value = "SynthetixNetwork"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
Set-TextValue $ws.Cells.Item(51, 4) "2.847"
$ws.Cells.Item(51, 5).Value = "  +33.28%  "
